$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108004
$ws.Cells.Item($row, 10).Value = "Papaya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 90
$ws.Cells.Item($row, 14).Value = 20000
$ws.Cells.Item($row, 15).Value = 20000
$ws.Cells.Item($row, 16).Value = 20000
$ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 2000
$ws.Cells.Item($row, 20).Value = 10
